$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Pin column (B) values for rows 2-7
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 6
$ws.Range("B7").Value = 7

# Fill in the remaining columns for row 7 (left_wrist)
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = "rotate_inside"
$ws.Range("E7").Value = 180
$ws.Range("F7").Value = "rotate_outside"

# Update the selected cell to F8
$ws.Range("F8").Select()
